# Add a second worksheet ("Sheet2") after the existing "Sheet1", populate it
# with the new shared-string text, size its first column, and move the
# "active/selected" tab state from Sheet1 to the newly added Sheet2 — mirroring
# what Excel does when a user adds a sheet, types a value, and tabs to it.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet directly after Sheet1 (Add() defaults to "before the
# active sheet", so pass the last sheet as the "After" anchor).
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# New content for the second sheet.
$ws2.Range("A1").Value = "Sheet 2 added new information for Azure Dev ops"

# After typing into A1 and hitting Enter, Excel's selection lands on A2.
[void]$ws2.Range("A2").Select()

# Match the "best fit" column width Excel would have computed for this text.
$ws2.Columns.Item(1).ColumnWidth = 45.16666666666666

Write-Host "Added Sheet2 with new content"
